# Daily attendance processing - 2025-09-28 15:16:18
# Updates "Recorded By", "Students" and "Status" for sessions that have now
# been recorded, and refreshes the dependent statistics tables.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Helper: write a percentage value as literal text without Excel's
# autodetection turning it into a numeric percent (which would change the
# cell's number format / style). We enter it with a leading apostrophe
# (forces text, just like typing it in the Excel UI) and then restore the
# original cell formatting by pasting formats only from an untouched donor
# cell that already carries the same style. ---
function Set-TextValue {
    param($range, [string]$value, $formatDonor)
    $range.Value = "'" + $value
    $formatDonor.Copy() | Out-Null
    $range.PasteSpecial(-4122) | Out-Null
}

# ----------------------------------------------------------------------
# Main attendance table: "Recorded By" (G), "Students" (H), "Status" (I)
# ----------------------------------------------------------------------

$ws.Range("G2").Value = "mennatulla.medhat@med.asu.edu.eg, servinaz@med.asu.edu.eg"
$ws.Range("H2").Value = "172/203"
$ws.Range("I2").Value = "Recorded"

$ws.Range("G9").Value = "Safa.hany@med.asu.edu.eg, Sara_nabil@med.asu.edu.eg, norhan.mohamed@med.asu.edu.eg, wessam.atef@med.asu.edu.eg"
$ws.Range("H9").Value = "168/203"
$ws.Range("I9").Value = "Recorded"

$ws.Range("G33").Value = "mennatulla.medhat@med.asu.edu.eg, servinaz@med.asu.edu.eg"
$ws.Range("H33").Value = "162/205"
$ws.Range("I33").Value = "Recorded"

$ws.Range("G40").Value = "Safa.hany@med.asu.edu.eg, Sara_nabil@med.asu.edu.eg, norhan.mohamed@med.asu.edu.eg, wessam.atef@med.asu.edu.eg"
$ws.Range("H40").Value = "164/205"
$ws.Range("I40").Value = "Recorded"

$ws.Range("G64").Value = "mennatulla.medhat@med.asu.edu.eg"
$ws.Range("H64").Value = "39/221"
$ws.Range("I64").Value = "Recorded"

$ws.Range("G65").Value = "mennatulla.medhat@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg"
$ws.Range("H65").Value = "63/221"
$ws.Range("I65").Value = "Recorded"

$ws.Range("G71").Value = "Omnia.Mohammed@med.asu.edu.eg, Safa.hany@med.asu.edu.eg, Sara_nabil@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg"
$ws.Range("H71").Value = "79/221"
$ws.Range("I71").Value = "Recorded"

$ws.Range("G95").Value = "Veronia.rafat@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, servinaz@med.asu.edu.eg"
$ws.Range("H95").Value = "2/132"
$ws.Range("I95").Value = "Recorded"

$ws.Range("G102").Value = "aml.awwad@med.asu.edu.eg, Safa.hany@med.asu.edu.eg, norhan.mohamed@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg"
$ws.Range("H102").Value = "28/132"
$ws.Range("I102").Value = "Recorded"

$ws.Range("G126").Value = "Veronia.rafat@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, servinaz@med.asu.edu.eg"
$ws.Range("H126").Value = "33/230"
$ws.Range("I126").Value = "Recorded"

$ws.Range("G133").Value = "aml.awwad@med.asu.edu.eg, Safa.hany@med.asu.edu.eg, norhan.mohamed@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg"
$ws.Range("H133").Value = "106/230"
$ws.Range("I133").Value = "Recorded"

$ws.Range("G157").Value = "servinaz@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg"
$ws.Range("H157").Value = "69/222"
$ws.Range("I157").Value = "Recorded"

$ws.Range("G158").Value = "mennatulla.medhat@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg"
$ws.Range("H158").Value = "59/222"
$ws.Range("I158").Value = "Recorded"

$ws.Range("G164").Value = "Omnia.Mohammed@med.asu.edu.eg, Safa.hany@med.asu.edu.eg, Sara_nabil@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg"
$ws.Range("H164").Value = "86/222"
$ws.Range("I164").Value = "Recorded"

# ----------------------------------------------------------------------
# Overall "Class Statistics" box (K3:L10)
# ----------------------------------------------------------------------

$ws.Range("L6").Value = 14    # Recorded Sessions
$ws.Range("L7").Value = 0     # Missing Sessions
Set-TextValue $ws.Range("L9") "7.5%" $ws.Range("L10")    # Coverage %
Set-TextValue $ws.Range("L10") "42.0%" $ws.Range("L9")   # Average Attendance %

# ----------------------------------------------------------------------
# "Group Statistics" box (K14:S20) - Recorded / Missing swap plus the
# recalculated Coverage % / Avg Attendance % columns.
# ----------------------------------------------------------------------

$ws.Range("O15").Value = 2
$ws.Range("P15").Value = 0
Set-TextValue $ws.Range("R15") "6.5%" $ws.Range("S16")
Set-TextValue $ws.Range("S15") "83.7%" $ws.Range("S16")

$ws.Range("O16").Value = 2
$ws.Range("P16").Value = 0
Set-TextValue $ws.Range("R16") "6.5%" $ws.Range("S17")
Set-TextValue $ws.Range("S16") "79.5%" $ws.Range("S17")

$ws.Range("O17").Value = 3
$ws.Range("P17").Value = 0
Set-TextValue $ws.Range("R17") "9.7%" $ws.Range("S18")
Set-TextValue $ws.Range("S17") "27.3%" $ws.Range("S18")

$ws.Range("O18").Value = 2
$ws.Range("P18").Value = 0
Set-TextValue $ws.Range("R18") "6.5%" $ws.Range("S19")
Set-TextValue $ws.Range("S18") "11.4%" $ws.Range("S19")

$ws.Range("O19").Value = 2
$ws.Range("P19").Value = 0
Set-TextValue $ws.Range("R19") "6.5%" $ws.Range("S20")
Set-TextValue $ws.Range("S19") "30.2%" $ws.Range("S20")

$ws.Range("O20").Value = 3
$ws.Range("P20").Value = 0
Set-TextValue $ws.Range("R20") "9.7%" $ws.Range("S15")
Set-TextValue $ws.Range("S20") "32.1%" $ws.Range("S15")

# ----------------------------------------------------------------------
# Column widths: widen "Recorded By" (G) and shrink "Status"-adjacent
# column I back down.
# ----------------------------------------------------------------------

$ws.Columns.Item(7).ColumnWidth = 50
$ws.Columns.Item(9).ColumnWidth = 10
